$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 21:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 524242
$ws.Range("C4").Value = 21366
$ws.Range("D4").Value = 28755
$ws.Range("E4").Value = 475264
$ws.Range("F4").Value = 11053
$ws.Range("G4").Value = 1476
$ws.Range("H4").Value = 20223

# Row 5 - España
$ws.Range("E5").Value = 86263
$ws.Range("G5").Value = 399
$ws.Range("H5").Value = 16480

# Row 8 - Alemania
$ws.Range("D8").Value = 57400
$ws.Range("E8").Value = 63742

# Row 19 - Austria
$ws.Range("B19").Value = 13799
$ws.Range("C19").Value = 239
$ws.Range("E19").Value = 6858

# Row 25 - India
$ws.Range("B25").Value = 8446
$ws.Range("C25").Value = 846
$ws.Range("E25").Value = 7189

# Row 29 - Noruega
$ws.Range("B29").Value = 6408
$ws.Range("C29").Value = 94
$ws.Range("E29").Value = 6258

# Row 36 - Pakistan
$ws.Range("E36").Value = 4163
$ws.Range("G36").Value = 20
$ws.Range("H36").Value = 86

# Row 88 - Costa Rica
$ws.Range("B88").Value = 577
$ws.Range("C88").Value = 19
$ws.Range("D88").Value = 49
$ws.Range("E88").Value = 525

# Row 114 - Isla de Man
$ws.Range("E114").Value = 112
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 2

# Row 163 - Siria
$ws.Range("E163").Value = 18
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 2
